$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.452.04'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.73%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.799.29'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.75%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.15%  '

# Row 6
$ws.Range('E6').Value = '  -0.06%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3796'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.52%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3447'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.83%  '

# Row 9
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.88'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.42%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.209'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.20%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07504'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.74%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9995'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.31%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +10.40%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.508'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.37%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.793.91'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.44%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.068'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.68%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001101'
$ws.Range('D17').Style = 'Normal'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06647'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.40%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '84.78'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.67%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.17%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.505'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.24%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.406.01'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.64%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.22%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.458'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.09%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.543'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +11.67%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.595'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +9.09%  '

# Row 28
$ws.Range('E28').Value = '  +11.73%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '150.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.54%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.996.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.52%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '133.67'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.15%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.049'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.53%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.176'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.87%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08669'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.73%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.39%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.687'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.90%  '

# Row 37
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6942'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.74%  '

# Row 38
$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.476'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.65%  '

# Row 39
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06379'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.17%  '

# Row 40
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2212'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.84%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.847'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.45%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02352'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.74%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.276'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.63%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.36%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6513'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.35%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.000'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.08%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.857'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.55%  '

# Row 48
$ws.Range('E48').Value = '  +5.04%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '130.48'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.59%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07197'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.61%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.71'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.97%  '
